# Adapt column header formatting to respective input file names.
# - Rename "_old" suffixed headers to "_FV2404"
# - Rename "_new" suffixed headers to "_FV2410"
# - Turn the header/data range into an actual Excel Table (ListObject)
# - Freeze the header row (pane split) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 21   # columns A..U
$lastRow = $ws.UsedRange.Rows.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Text
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2404")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2410")
    }
}

# Determine the full data range (A1:U66)
$startCell = $ws.Cells.Item(1, 1)
$endCell = $ws.Cells.Item($lastRow, $lastCol)
$tableRange = $ws.Range($startCell, $endCell)

# Create an Excel Table (ListObject) over the header + data range
$listObject = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# Freeze the header row: select the first cell below the header, then
# turn on FreezePanes so Excel records a plain (non-split) frozen pane.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
